# Repull data, push all data, mean calculation
# Updates the dSF (column F) values on Sheet1 for the rows whose
# underlying data changed after the repull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    7  = -2
    14 = 3
    18 = 1
    20 = 5
    21 = -2
    24 = -6
    26 = -2
    27 = 3
    34 = -3
    35 = -2
    37 = -4
    39 = -12
    41 = -8
    43 = 4
    45 = -3
    47 = -4
    49 = -3
    57 = 1
    61 = -3
    62 = 5
    65 = 5
    66 = 0
    70 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
